$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to a text format before writing so that
# Excel does not auto-coerce numeric-looking strings (e.g. "148.52") into
# floating point numbers. This mirrors the source data which stores these
# as plain text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.573.09"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.594.91"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "609.37"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "148.52"
$ws.Range("E6").Value = "  +2.51%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("B9").Value = "Toncoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D9").Value = "8.05"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "4.209.75"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "0.0000209"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "29.81"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "3.585.79"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "66.691.53"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "11.51"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "15.09"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "427.06"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "0.617"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "78.81"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "3.743.02"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").Value = "9.36"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "3.595.04"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "0.158"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "25.44"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "7.85"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "5.64"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "177.78"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").Value = "0.0856"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "0.898"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "2.56"
$ws.Range("E44").Value = "  +8.21%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "1.18"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "25.03"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").Value = "24.09"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "7.18"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "0.953"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "2.419.55"
$ws.Range("E51").Value = "  +4.81%  "

# Restore default (no explicit) cell formatting so the workbook styles stay
# equivalent to the original, while keeping the text values intact.
$ws.Range("D2:E51").ClearFormats()
